$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '62.906.16'
$ws.Range("E2").Value = '  -5.27%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.106.47'
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '558.10'
$ws.Range("E5").Value = '  -4.86%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '162.15'
$ws.Range("E6").Value = '  -10.02%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.586'
$ws.Range("E8").Value = '  -7.77%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3.100.59'
$ws.Range("E9").Value = '  -5.73%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.73'
$ws.Range("E10").Value = '  -1.74%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.115'
$ws.Range("E11").Value = '  -7.91%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.377'
$ws.Range("E12").Value = '  -6.09%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.642.40'
$ws.Range("E13").Value = '  -5.82%  '
$ws.Range("E14").Value = '  -1.87%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '62.994.72'
$ws.Range("E15").Value = '  -5.12%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '24.53'
$ws.Range("E16").Value = '  -7.56%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.103.60'
$ws.Range("E17").Value = '  -7.22%  '
$ws.Range("E18").Value = '  -6.65%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '400.81'
$ws.Range("E19").Value = '  -6.08%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.34'
$ws.Range("E20").Value = '  -4.97%  '
$ws.Range("E21").Value = '  -5.58%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.03'
$ws.Range("E22").Value = '  -3.53%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.999'
$ws.Range("E23").Value = '  +0.09%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.70'
$ws.Range("E24").Value = '  -0.33%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '67.49'
$ws.Range("E25").Value = '  -5.69%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.198'
$ws.Range("E26").Value = '  -3.76%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.475'
$ws.Range("E27").Value = '  -7.30%  '
$ws.Range("E28").Value = '  -11.85%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.61'
$ws.Range("E29").Value = '  -5.43%  '
$ws.Range("E30").Value = '  +0.05%  '
$ws.Range("E32").Value = '  -7.28%  '
$ws.Range("E33").Value = '  -6.32%  '
$ws.Range("E34").Value = '  -6.44%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.16'
$ws.Range("E35").Value = '  -5.82%  '
$ws.Range("E36").Value = '  -7.07%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '151.81'
$ws.Range("E37").Value = '  -4.72%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.31'
$ws.Range("E38").Value = '  -7.81%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.698.70'
$ws.Range("E39").Value = '  -5.90%  '
$ws.Range("E40").Value = '  -8.28%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '23.37'
$ws.Range("E41").Value = '  -10.61%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.00'
$ws.Range("E42").Value = '  -6.92%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '38.15'
$ws.Range("E43").Value = '  -3.98%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.691'
$ws.Range("E44").Value = '  -7.85%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0601'
$ws.Range("E45").Value = '  -8.16%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0254'
$ws.Range("E46").Value = '  -6.28%  '
$ws.Range("E47").Value = '  -12.40%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '282.64'
$ws.Range("E48").Value = '  -9.19%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '20.68'
$ws.Range("E49").Value = '  -8.49%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.00'
$ws.Range("E50").Value = '  +0.00%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0971'
$ws.Range("E51").Value = '  -5.30%  '
